# Update TPM-derived NATMI Ligand-Receptor expression/specificity metrics
# for the Pdgfb-Lrp1 sheet with freshly recomputed values (new TPM input).
#
# Columns G/H/I/J  -> Ligand average/total expression value + derived
#                     specificity (depends only on the Sending cluster).
# Columns M/N/O/P  -> Receptor average/total expression value + derived
#                     specificity (depends only on the Target cluster).
# Columns Q/R/S/T  -> Edge average/total expression weight + derived
#                     specificity (depends on the Sending/Target pair).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 45.01222466666667
$ws.Range("H2").Value = 135.036674
$ws.Range("I2").Value = 0.7482903203664146
$ws.Range("J2").Value = 0.7482903203664146
$ws.Range("M2").Value = 1.918906333333333
$ws.Range("N2").Value = 5.756718999999999
$ws.Range("O2").Value = 0.006524019162508824
$ws.Range("P2").Value = 0.006524019162508824
$ws.Range("Q2").Value = 86.37424299028955
$ws.Range("R2").Value = 777.368186912606
$ws.Range("S2").Value = 0.004881860389190356
$ws.Range("T2").Value = 0.004881860389190356
$ws.Range("G3").Value = 45.01222466666667
$ws.Range("H3").Value = 135.036674
$ws.Range("I3").Value = 0.7482903203664146
$ws.Range("J3").Value = 0.7482903203664146
$ws.Range("O3").Value = 0.6163557430885885
$ws.Range("P3").Value = 0.6163557430885885
$ws.Range("Q3").Value = 8160.193800154588
$ws.Range("R3").Value = 73441.74420139128
$ws.Range("S3").Value = 0.4612130364554395
$ws.Range("T3").Value = 0.4612130364554395
$ws.Range("G4").Value = 45.01222466666667
$ws.Range("H4").Value = 135.036674
$ws.Range("I4").Value = 0.7482903203664146
$ws.Range("J4").Value = 0.7482903203664146
$ws.Range("M4").Value = 29.04767233333333
$ws.Range("N4").Value = 87.143017
$ws.Range("O4").Value = 0.09875811426384234
$ws.Range("P4").Value = 0.09875811426384236
$ws.Range("Q4").Value = 1307.500353111718
$ws.Range("R4").Value = 11767.50317800546
$ws.Range("S4").Value = 0.07389974096127357
$ws.Range("T4").Value = 0.07389974096127358
$ws.Range("G5").Value = 45.01222466666667
$ws.Range("H5").Value = 135.036674
$ws.Range("I5").Value = 0.7482903203664146
$ws.Range("J5").Value = 0.7482903203664146
$ws.Range("M5").Value = 81.87450533333333
$ws.Range("N5").Value = 245.623516
$ws.Range("O5").Value = 0.2783621234850603
$ws.Range("P5").Value = 0.2783621234850603
$ws.Range("Q5").Value = 3685.353628536198
$ws.Range("R5").Value = 33168.18265682578
$ws.Range("S5").Value = 0.2082956825605112
$ws.Range("T5").Value = 0.2082956825605112
$ws.Range("G6").Value = 2.766295666666667
$ws.Range("H6").Value = 8.298887000000001
$ws.Range("I6").Value = 0.04598733535094824
$ws.Range("J6").Value = 0.04598733535094825
$ws.Range("M6").Value = 1.918906333333333
$ws.Range("N6").Value = 5.756718999999999
$ws.Range("O6").Value = 0.006524019162508824
$ws.Range("P6").Value = 0.006524019162508824
$ws.Range("Q6").Value = 5.308262274639222
$ws.Range("R6").Value = 47.774360471753
$ws.Range("S6").Value = 0.0003000222570623058
$ws.Range("T6").Value = 0.0003000222570623058
$ws.Range("G7").Value = 2.766295666666667
$ws.Range("H7").Value = 8.298887000000001
$ws.Range("I7").Value = 0.04598733535094824
$ws.Range("J7").Value = 0.04598733535094825
$ws.Range("O7").Value = 0.6163557430885885
$ws.Range("P7").Value = 0.6163557430885885
$ws.Range("Q7").Value = 501.497291362371
$ws.Range("R7").Value = 4513.475622261339
$ws.Range("S7").Value = 0.02834455825289782
$ws.Range("T7").Value = 0.02834455825289782
$ws.Range("G8").Value = 2.766295666666667
$ws.Range("H8").Value = 8.298887000000001
$ws.Range("I8").Value = 0.04598733535094824
$ws.Range("J8").Value = 0.04598733535094825
$ws.Range("M8").Value = 29.04767233333333
$ws.Range("N8").Value = 87.143017
$ws.Range("O8").Value = 0.09875811426384234
$ws.Range("P8").Value = 0.09875811426384236
$ws.Range("Q8").Value = 80.35445010245323
$ws.Range("R8").Value = 723.1900509220791
$ws.Range("S8").Value = 0.004541622519278582
$ws.Range("T8").Value = 0.004541622519278584
$ws.Range("G9").Value = 2.766295666666667
$ws.Range("H9").Value = 8.298887000000001
$ws.Range("I9").Value = 0.04598733535094824
$ws.Range("J9").Value = 0.04598733535094825
$ws.Range("M9").Value = 81.87450533333333
$ws.Range("N9").Value = 245.623516
$ws.Range("O9").Value = 0.2783621234850603
$ws.Range("P9").Value = 0.2783621234850603
$ws.Range("Q9").Value = 226.4890893140769
$ws.Range("R9").Value = 2038.401803826692
$ws.Range("S9").Value = 0.01280113232170953
$ws.Range("T9").Value = 0.01280113232170953
$ws.Range("G10").Value = 12.37490333333333
$ws.Range("H10").Value = 37.12471
$ws.Range("I10").Value = 0.2057223442826371
$ws.Range("J10").Value = 0.2057223442826371
$ws.Range("M10").Value = 1.918906333333333
$ws.Range("N10").Value = 5.756718999999999
$ws.Range("O10").Value = 0.006524019162508824
$ws.Range("P10").Value = 0.006524019162508824
$ws.Range("Q10").Value = 23.74628038072111
$ws.Range("R10").Value = 213.71652342649
$ws.Range("S10").Value = 0.001342136516256162
$ws.Range("T10").Value = 0.001342136516256162
$ws.Range("G11").Value = 12.37490333333333
$ws.Range("H11").Value = 37.12471
$ws.Range("I11").Value = 0.2057223442826371
$ws.Range("J11").Value = 0.2057223442826371
$ws.Range("O11").Value = 0.6163557430885885
$ws.Range("P11").Value = 0.6163557430885885
$ws.Range("Q11").Value = 2243.426318205505
$ws.Range("R11").Value = 20190.83686384954
$ws.Range("S11").Value = 0.1267981483802512
$ws.Range("T11").Value = 0.1267981483802512
$ws.Range("G12").Value = 12.37490333333333
$ws.Range("H12").Value = 37.12471
$ws.Range("I12").Value = 0.2057223442826371
$ws.Range("J12").Value = 0.2057223442826371
$ws.Range("M12").Value = 29.04767233333333
$ws.Range("N12").Value = 87.143017
$ws.Range("O12").Value = 0.09875811426384234
$ws.Range("P12").Value = 0.09875811426384236
$ws.Range("Q12").Value = 359.4621371833412
$ws.Range("R12").Value = 3235.15923465007
$ws.Range("S12").Value = 0.02031675078329019
$ws.Range("T12").Value = 0.02031675078329019
$ws.Range("G13").Value = 12.37490333333333
$ws.Range("H13").Value = 37.12471
$ws.Range("I13").Value = 0.2057223442826371
$ws.Range("J13").Value = 0.2057223442826371
$ws.Range("M13").Value = 81.87450533333333
$ws.Range("N13").Value = 245.623516
$ws.Range("O13").Value = 0.2783621234850603
$ws.Range("P13").Value = 0.2783621234850603
$ws.Range("Q13").Value = 1013.189088964484
$ws.Range("R13").Value = 9118.701800680359
$ws.Range("S13").Value = 0.05726530860283952
$ws.Range("T13").Value = 0.05726530860283952
